# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" table: a handful of countries (Rusia,
# Indonesia, Lituania, Kazajistan) moved up in the ranking (column B,
# "Casos totales") with newly updated figures, pushing the countries
# that used to occupy those ranks down by one row with their previous
# figures carried along. Also bump the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 28 de Marzo de 2020 a las 09:59"

# Row 20
$ws.Range("B20").Value = 3773
$ws.Range("C20").Value = 2
$ws.Range("E20").Value = 3747

# Row 36
$ws.Range("A36").Value = "Rusia"
$ws.Range("B36").Value = 1264
$ws.Range("C36").Value = 228
$ws.Range("D36").Value = 49
$ws.Range("E36").Value = 1211
$ws.Range("F36").Value = 8
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 4

# Row 37
$ws.Range("A37").Value = "Tailandia"
$ws.Range("B37").Value = 1245
$ws.Range("C37").Value = 109
$ws.Range("D37").Value = 97
$ws.Range("E37").Value = 1142
$ws.Range("F37").Value = 11
$ws.Range("G37").Value = 1
$ws.Range("H37").Value = 6

# Row 38
$ws.Range("A38").Value = "Sudafrica"
$ws.Range("B38").Value = 1170
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 31
$ws.Range("E38").Value = 1138
$ws.Range("F38").Value = 7
$ws.Range("H38").Value = 1

# Row 39
$ws.Range("A39").Value = "Finlandia"
$ws.Range("B39").Value = 1165
$ws.Range("C39").Value = 124
$ws.Range("D39").Value = 10
$ws.Range("E39").Value = 1148
$ws.Range("F39").Value = 32
$ws.Range("H39").Value = 7

# Row 40
$ws.Range("B40").Value = 1155
$ws.Range("C40").Value = 109
$ws.Range("D40").Value = 59
$ws.Range("E40").Value = 994
$ws.Range("G40").Value = 15
$ws.Range("H40").Value = 102

# Row 41
$ws.Range("A41").Value = "Arabia Saudita"
$ws.Range("B41").Value = 1104
$ws.Range("D41").Value = 35
$ws.Range("E41").Value = 1066
$ws.Range("F41").Value = 6
$ws.Range("H41").Value = 3

# Row 43
$ws.Range("B43").Value = 908
$ws.Range("C43").Value = 21
$ws.Range("E43").Value = 805

# Row 67
$ws.Range("A67").Value = "Lituania"
$ws.Range("B67").Value = 382
$ws.Range("C67").Value = 24
$ws.Range("D67").Value = 1
$ws.Range("E67").Value = 376
$ws.Range("F67").Value = 2
$ws.Range("H67").Value = 5

# Row 68
$ws.Range("A68").Value = "Armenia"
$ws.Range("B68").Value = 372
$ws.Range("C68").Value = 43
$ws.Range("D68").Value = 28
$ws.Range("E68").Value = 343
$ws.Range("F68").Value = 6
$ws.Range("H68").Value = 1

# Row 69
$ws.Range("B69").Value = 358
$ws.Range("C69").Value = 13
$ws.Range("E69").Value = 324

# Row 73
$ws.Range("B73").Value = 305
$ws.Range("C73").Value = 25
$ws.Range("E73").Value = 304

# Row 85
$ws.Range("A85").Value = "Kazajistan"
$ws.Range("B85").Value = 204
$ws.Range("C85").Value = 54
$ws.Range("D85").Value = 14
$ws.Range("E85").Value = 189
$ws.Range("F85").Value = 0
$ws.Range("H85").Value = 1

# Row 86
$ws.Range("A86").Value = "Moldavia"
$ws.Range("B86").Value = 199
$ws.Range("D86").Value = 2
$ws.Range("E86").Value = 195
$ws.Range("F86").Value = 33
$ws.Range("H86").Value = 2

# Row 87
$ws.Range("A87").Value = "Albania"
$ws.Range("B87").Value = 186
$ws.Range("D87").Value = 31
$ws.Range("E87").Value = 147
$ws.Range("F87").Value = 3
$ws.Range("H87").Value = 8

# Row 88
$ws.Range("A88").Value = "Burkina Faso"
$ws.Range("B88").Value = 180
$ws.Range("C88").Value = 0
$ws.Range("D88").Value = 12
$ws.Range("E88").Value = 159
$ws.Range("F88").Value = 0
$ws.Range("H88").Value = 9

# Row 89
$ws.Range("A89").Value = "Vietnam"
$ws.Range("B89").Value = 169
$ws.Range("C89").Value = 6
$ws.Range("D89").Value = 20
$ws.Range("E89").Value = 149
$ws.Range("F89").Value = 3
$ws.Range("H89").Value = 0

# Row 90
$ws.Range("A90").Value = "Azerbaiyan"
$ws.Range("B90").Value = 165
$ws.Range("E90").Value = 147
$ws.Range("F90").Value = 6
$ws.Range("H90").Value = 3

# Row 91
$ws.Range("A91").Value = "Republica de Chipre"
$ws.Range("B91").Value = 162
$ws.Range("C91").Value = 0
$ws.Range("D91").Value = 15
$ws.Range("E91").Value = 142
$ws.Range("F91").Value = 3
$ws.Range("H91").Value = 5

# Row 92
$ws.Range("A92").Value = "Oman"
$ws.Range("B92").Value = 152
$ws.Range("C92").Value = 21
$ws.Range("D92").Value = 23
$ws.Range("E92").Value = 129
$ws.Range("H92").Value = 0
